$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VSS")

# Rows 10 (VSS-09) and 11 (VSS-10) describe the SMB backup/restore scenarios
# ("Backup/restore using VHDx/VHD on SMB"), but their embedded LISA test XML
# (column G) had the <testName> copy-pasted from the CSV testcases above.
# Correct the test names so they match the SMB scenario they actually cover.
$g10 = $ws.Range("G10").Value2
$ws.Range("G10").Value = $g10.Replace("VSS_BackupRestore_CSV_VHDX", "VSS_BackupRestore_SMB_VHDX")

$g11 = $ws.Range("G11").Value2
$ws.Range("G11").Value = $g11.Replace("VSS_BackupRestore_CSV_VHD", "VSS_BackupRestore_SMB_VHD")

# Restore the viewport/selection state recorded when the sheet was last saved.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G11").Select()
